$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (methods): split "API, data wrangling, maskinlæring, tekstklassifisering"
# into separate cells B5:E5
$ws.Range("B5").Value = "API"
$ws.Range("C5").Value = "data wrangling"
$ws.Range("D5").Value = "maskinlæring"
$ws.Range("E5").Value = "tekstklassifisering"

# Row 6 (themes): split "Sosiale medier, tekstanalyse, misinformasjon"
# into separate cells B6:D6
$ws.Range("B6").Value = "Sosiale medier"
$ws.Range("C6").Value = "tekstanalyse"
$ws.Range("D6").Value = "misinformasjon"

# Update selection to match the recorded state (B7 selected)
$ws.Range("B7").Select()
